# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Temperature sheets,
# mirroring the "Auto update Excel log" automated appender.
#
# Columns (all sheets): A=Date, B=Timestamp, C=Hour, D=Location, E=Value, F=Status
#
# Column A ("Date", e.g. "2026-02-06") and, on the Humidity sheet, column E
# ("Value", e.g. "68.0%") look like a date / a percentage to Excel's smart
# entry, so NumberFormat is forced to Text ("@") for those ranges before the
# values are written - otherwise Excel would silently coerce them into a
# date serial number / a percentage number instead of keeping the literal
# log text.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [string]$Csv,
        [bool]$ValueColumnIsPercent
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $lines = $Csv.Split("`n")
    $rowCount = 0
    foreach ($line in $lines) {
        $trimmed = $line.Trim()
        if ($trimmed.Length -gt 0) {
            $rowCount = $rowCount + 1
        }
    }
    $endRow = $StartRow + $rowCount - 1

    # Force the Date column (A) to Text so "2026-02-06" isn't turned into a
    # date serial number.
    $ws.Range("A" + $StartRow + ":A" + $endRow).NumberFormat = "@"
    if ($ValueColumnIsPercent) {
        # Humidity's Value column holds literal strings like "68.0%" which
        # Excel would otherwise coerce into a percentage number.
        $ws.Range("E" + $StartRow + ":E" + $endRow).NumberFormat = "@"
    }

    $r = $StartRow
    foreach ($line in $lines) {
        $trimmed = $line.Trim()
        if ($trimmed.Length -gt 0) {
            $fields = $trimmed.Split(",")
            $ws.Cells.Item($r, 1).Value = $fields[0]
            $ws.Cells.Item($r, 2).Value = $fields[1]
            $ws.Cells.Item($r, 3).Value = $fields[2]
            $ws.Cells.Item($r, 4).Value = $fields[3]
            $ws.Cells.Item($r, 5).Value = $fields[4]
            $ws.Cells.Item($r, 6).Value = $fields[5]
            $r = $r + 1
        }
    }

    # Re-apply the plain "Normal" cell style now that the literal text is
    # safely stored - this drops the temporary Text number-format override
    # again (matching the source log rows, which carry no explicit style)
    # while leaving the cells' stored values as literal text.
    $ws.Range("A" + $StartRow + ":A" + $endRow).Style = "Normal"
    if ($ValueColumnIsPercent) {
        $ws.Range("E" + $StartRow + ":E" + $endRow).Style = "Normal"
    }
}

# --- PIR sheet: rows 540-552 (A1:F539 -> A1:F552) ---
$pirCsv = @"
2026-02-06,10:24:48,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:24:52,10:00,Bathroom,Motion Detected,Active
2026-02-06,10:24:59,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:01,10:00,Bathroom,Motion Detected,Active
2026-02-06,10:25:06,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:11,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:16,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:20,10:00,Bathroom,Motion Detected,Active
2026-02-06,10:25:25,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:30,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:35,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:40,10:00,Bathroom,No Motion,Inactive
2026-02-06,10:25:45,10:00,Bathroom,No Motion,Inactive
"@
Add-LogRows "PIR" 540 $pirCsv $false

# --- Humidity sheet: rows 379-390 (A1:F378 -> A1:F390) ---
$humidityCsv = @"
2026-02-06,10:24:50,10:00,Bathroom,68.0%,Active
2026-02-06,10:24:53,10:00,Bathroom,67.1%,Active
2026-02-06,10:24:57,10:00,Bathroom,67.9%,Active
2026-02-06,10:25:02,10:00,Bathroom,67.0%,Active
2026-02-06,10:25:07,10:00,Bathroom,68.0%,Active
2026-02-06,10:25:12,10:00,Bathroom,67.9%,Active
2026-02-06,10:25:17,10:00,Bathroom,68.1%,Active
2026-02-06,10:25:22,10:00,Bathroom,68.1%,Active
2026-02-06,10:25:27,10:00,Bathroom,68.2%,Active
2026-02-06,10:25:37,10:00,Bathroom,68.1%,Active
2026-02-06,10:25:42,10:00,Bathroom,67.0%,Active
2026-02-06,10:25:47,10:00,Bathroom,67.8%,Active
"@
Add-LogRows "Humidity" 379 $humidityCsv $true

# --- Temperature sheet: rows 379-389 (A1:F378 -> A1:F389) ---
$temperatureCsv = @"
2026-02-06,10:24:51,10:00,Bathroom,28.4C,Active
2026-02-06,10:24:54,10:00,Bathroom,28.4C,Active
2026-02-06,10:24:58,10:00,Bathroom,28.3C,Active
2026-02-06,10:25:03,10:00,Bathroom,28.3C,Active
2026-02-06,10:25:08,10:00,Bathroom,28.4C,Active
2026-02-06,10:25:13,10:00,Bathroom,28.3C,Active
2026-02-06,10:25:19,10:00,Bathroom,28.4C,Active
2026-02-06,10:25:23,10:00,Bathroom,28.4C,Active
2026-02-06,10:25:28,10:00,Bathroom,28.5C,Active
2026-02-06,10:25:39,10:00,Bathroom,28.4C,Active
2026-02-06,10:25:44,10:00,Bathroom,28.4C,Active
"@
Add-LogRows "Temperature" 379 $temperatureCsv $false
